# Add a new column C of image-filename labels to Sheet1.
# The values are written in the exact order the author originally entered
# them so that the shared-string table gets appended in the same sequence
# as the source workbook (row 12 before row 11 - see below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value  = "thienyetcugiai.jpg"
$ws.Range("C4").Value  = "12-cung-hoang-dao.jpg"
$ws.Range("C7").Value  = "su-nghiep-cung-thien-binh.jpg"
$ws.Range("C8").Value  = "virgo.jpg"
$ws.Range("C9").Value  = "cung-su-tu-13.jpg"
$ws.Range("C10").Value = "12-cung-hoang-dao2.jpg"
$ws.Range("C12").Value = "kim-nguu-va-ma-ket.jpg"
$ws.Range("C11").Value = "12-cung-hoang-dao3.jpg"
$ws.Range("C14").Value = "thien-binh-va-thien-yet-ong-nam-lang-ba-nua-can"
$ws.Range("C23").Value = "sutu-25e16.jpg"
$ws.Range("C26").Value = "bachduongnhanma.jpg"
$ws.Range("C38").Value = "baobinh.jpg"
$ws.Range("C40").Value = "maketbaobinh.jpg"
$ws.Range("C52").Value = "thienbinhvabaobinh.jpg"
$ws.Range("C54").Value = "thienbinhnhanma.jpg"
$ws.Range("C76").Value = "thienyetcugiai2.jpg"
$ws.Range("C92").Value = "kim-nguu-va-ma-ket2.jpg"

# New column C width (matches the author's "best fit" auto-resize as closely
# as this host's pixel-snapped ColumnWidth setter allows).
$ws.Columns.Item(3).ColumnWidth = 43

# Selection / view tweaks the author made while editing.
$ws.Range("E84").Select()
